$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.205.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.655.38'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  +0.72%  '
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.56'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.886.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.654.90'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.189.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0737'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.65%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.39%  '
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0514'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  +2.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.266.85'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -1.61%  '
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.826'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.808'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.796.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0517'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0976'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('E51').Value = '  +0.20%  '
